# "Actualizacion desde MV -datos-"
# Append 3 new daily rows (04-10-2021, 05-10-2021, 06-10-2021) to the bottom
# of the "Facilidades permanentes" table on Sheet1, right after the existing
# last row (190, date 01-10-2021).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 190
$newRows = @(
    @{ Date = "04-10-2021"; FPD = 5726450; FPL = 81401 },
    @{ Date = "05-10-2021"; FPD = 6195650; FPL = 8883 },
    @{ Date = "06-10-2021"; FPD = 6223850; FPL = 0 }
)

# Scratch area (well below the used range) used to build the date text via a
# formula. A formula result of type string is stored as plain text (shared
# string) with no special "looks like a date" auto-formatting/style applied,
# unlike assigning the literal text straight to .Value which Excel would
# otherwise auto-convert into a date serial number with a date number format.
$scratchRow = $lastRow + 200

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $scratchCell = $ws.Cells.Item($scratchRow + $i, 1)
    $scratchCell.Formula = '="' + $newRows[$i].Date + '"'
}

$scratchRange = $ws.Range($ws.Cells.Item($scratchRow, 1), $ws.Cells.Item($scratchRow + $newRows.Count - 1, 1))
$scratchRange.Copy()
$destRange = $ws.Range($ws.Cells.Item($lastRow + 1, 1), $ws.Cells.Item($lastRow + $newRows.Count, 1))
$destRange.PasteSpecial(-4163)

# Remove the scratch rows again, leaving no trace (no leftover formulas, no
# leftover styles).
$scratchRange.EntireRow.Delete()

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $lastRow + 1 + $i
    $ws.Cells.Item($r, 2).Value = $newRows[$i].FPD
    $ws.Cells.Item($r, 3).Value = $newRows[$i].FPL
}
